# Added filtering options for the Component Analysis
#
# The evaluation window rolled forward by one period: a new evaluation
# (previously absent) is inserted at the top of the data block (row 2),
# every existing data row's statistics shift down by one row, and the
# oldest data row (old row 11) falls out of the window. Column A (the
# Q0..Q9 labels) stays positional/unchanged - only the numeric statistics
# in columns B:G move.
#
# NOTE: this engine's Range.Value/.Value2 setter only reliably applies for
# single-cell assignments - assigning a PowerShell array to a multi-cell
# range silently/erroneously no-ops - so every cell is written individually
# via Cells.Item(row, col).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 11
$firstCol = 2   # column B
$lastCol = 7    # column G

# Snapshot the current (pre-edit) per-cell values for the data block so the
# down-shift below can't clobber a value before it has been read.
$old = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $old["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Shift every existing row's statistics down by one row: old row r's values
# move into row r+1 (processed bottom-up so row r hasn't been overwritten
# yet when it's read for row r+1).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $old["$($r - 1),$c"]
    }
}

# Write the brand-new evaluation's statistics into the now-vacated top row.
$newRow = @(0.08266386729847572, 1.266710845429791, 10.58350078540567, 3.253229285710688, 3.288115649630924, 46)
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item($firstRow, $c).Value2 = $newRow[$c - $firstCol]
}
